$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.1070994234211895,
    -0.0099999994805806125,
    -0.008999999487731003,
    0.28399077636663961,
    -0.0059999995050104005,
    -0.005999999489596064,
    -0.01999999939887509,
    -0.019999999396590695,
    -0.0059999994850095106,
    -0.0059999994845014726,
    -0.0044999994940830845,
    -0.00599999948469776,
    0.02037638542821707,
    -0.011999999447794174,
    -0.0059999994864297079,
    0.025420419632430136,
    -0.0059999994850121752,
    -0.0089999994652769644,
    -0.048302720236149899,
    -0.008999999472264264,
    -0.0089999994694305308,
    -0.051816882489932325,
    -0.0089999994836951203,
    -0.041999999264207943,
    -0.041999999260199594,
    -0.0041681432169369259,
    -0.0059999994848345395,
    -0.0059999994734036832,
    -0.011999999426913988,
    -0.019999999371881572,
    -0.014999999398547459,
    -0.020999999359218258,
    -0.0059999994550770097
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
